# Append the latest daily portfolio snapshot (2025-09-21) as a new row
# right after the existing data, extending the table from A1:D36 to A1:D37.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37

# Column A stores plain text dates (e.g. "2025-08-17", "2025-09-20") rather
# than real date values elsewhere in this sheet, so use a leading apostrophe
# to force this literal to stay text instead of being auto-converted to a
# date serial by Excel.
$ws.Cells.Item($row, 1).Value = "'2025-09-21"
$ws.Cells.Item($row, 2).Value = 60.40000152587891
$ws.Cells.Item($row, 3).Value = 707.4500122070312
$ws.Cells.Item($row, 4).Value = 336.5499877929688
